$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New section appended below existing content: "Rules for Nouns and its Plurals".
# Cell writes are ordered to reproduce the original shared-string build order
# (note rows 149 vs 151: "Potato"/"Potatoes" were interned before the "Rule 6"
# heading text, even though the heading sits on an earlier row).

$ws.Range("A117").Value = 'Rules for Nouns and its Plurals'
$ws.Range("A118").Value = 'Rule 1 : Add an (s) to form the plural of most of words. '
$ws.Range("A119").Value = 'Ex :'
$ws.Range("A132").Value = 'Ex :'
$ws.Range("A137").Value = 'Ex :'
$ws.Range("A150").Value = 'Ex :'
$ws.Range("A120").Value = 'Book'
$ws.Range("B120").Value = 'Books'
$ws.Range("A121").Value = 'Bug'
$ws.Range("B121").Value = 'Bugs'
$ws.Range("A122").Value = 'Friend'
$ws.Range("B122").Value = 'Friends'
$ws.Range("A124").Value = 'Rule 2 :For world that ends in (s,z,x,ch,sh,ss) add (es) to form the plural.'
$ws.Range("A125").Value = 'Ex:'
$ws.Range("A126").Value = 'Quiz'
$ws.Range("B126").Value = 'Quizes'
$ws.Range("A127").Value = 'Church'
$ws.Range("B127").Value = 'Churches'
$ws.Range("A128").Value = 'Box'
$ws.Range("B128").Value = 'Boxes'
$ws.Range("A129").Value = 'Class'
$ws.Range("B129").Value = 'Classes'
$ws.Range("A131").Value = 'Rule 3 : If the word ends in a vowel + (y) [ay,ey,iy,oy,uy] add an s to to word.'
$ws.Range("A133").Value = 'Boy'
$ws.Range("B133").Value = 'Boys'
$ws.Range("A134").Value = 'Stay'
$ws.Range("B134").Value = 'Stays'
$ws.Range("A136").Value = 'Rule 4 : If the word ends in a consonant + (y) change the (y) into (i) and add (es) to form the plural.'
$ws.Range("A138").Value = 'Enemy'
$ws.Range("B138").Value = 'Enimies'
$ws.Range("A139").Value = 'Cherry'
$ws.Range("B139").Value = 'Cherries'
$ws.Range("A140").Value = 'Baby'
$ws.Range("B140").Value = 'Babies'
$ws.Range("A142").Value = 'Rule 5 : For words that end in (f), drop the (f) change it to (v) and add (es) for the pluerl.'
$ws.Range("A143").Value = 'For words that end in (fe), change the (f) to (v) and add (es) to form the plural.'
$ws.Range("A144").Value = 'Eg :'
$ws.Range("A145").Value = 'Half'
$ws.Range("B145").Value = 'Halves'
$ws.Range("A146").Value = 'Wife'
$ws.Range("B146").Value = 'Wives'
$ws.Range("A147").Value = 'Wolf'
$ws.Range("B147").Value = 'Wolves'
$ws.Range("A151").Value = 'Potato'
$ws.Range("B151").Value = 'Potatoes'
$ws.Range("A149").Value = 'Rule 6: For nouns that ends in consonents + (o) add (es) to make the noun plural.'
$ws.Range("A152").Value = 'Hero'
$ws.Range("B152").Value = 'Heroes'
$ws.Range("A153").Value = 'Volcano'
$ws.Range("B153").Value = 'Volcanoes'
$ws.Range("A155").Value = 'Irregular Plural'
$ws.Range("A156").Value = 'Some nouns do not change into plural.'
$ws.Range("A157").Value = 'One deer'
$ws.Range("B157").Value = 'Two deer'
$ws.Range("A158").Value = 'one sheep'
$ws.Range("B158").Value = 'Two sheep'
$ws.Range("A160").Value = 'Some nouns are always plural'
$ws.Range("A161").Value = 'Glases'
$ws.Range("A162").Value = 'Pyjamas'
$ws.Range("A163").Value = 'Shorts'
$ws.Range("A164").Value = 'Cloths'

# Section heading style (bold), matching the other section headers in the sheet.
$ws.Range("A155").Font.Bold = $true

# Leave selection where the user ended up after typing the new content.
$ws.Range("A165").Select()
$excel.ActiveWindow.ScrollRow = 134
$excel.ActiveWindow.ScrollColumn = 1
